$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 4 trailing rows (table shrinks from 72 to 68 rows)
$ws.Range("A69:A72").EntireRow.Delete() | Out-Null

# Force columns C:H to text so numeric-looking / date-looking strings are preserved as text
$ws.Range("C2:H68").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = '?'
$ws.Cells.Item(2, 2).Value = 489
$ws.Cells.Item(2, 3).Value = '15440000'
$ws.Cells.Item(2, 4).Value = '17880828'
$ws.Cells.Item(2, 5).Value = '1544-07-02'
$ws.Cells.Item(2, 6).Value = '1788-08-28'
$ws.Cells.Item(2, 7).Value = '1544'
$ws.Cells.Item(2, 8).Value = '1788'
$ws.Cells.Item(3, 1).Value = 'Coimbra'
$ws.Cells.Item(3, 2).Value = 57
$ws.Cells.Item(3, 3).Value = '15420000'
$ws.Cells.Item(3, 4).Value = '17420427'
$ws.Cells.Item(3, 5).Value = '1542-07-02'
$ws.Cells.Item(3, 6).Value = '1742-04-27'
$ws.Cells.Item(3, 7).Value = '1542'
$ws.Cells.Item(3, 8).Value = '1742'
$ws.Cells.Item(4, 1).Value = 'Lisboa'
$ws.Cells.Item(4, 2).Value = 45
$ws.Cells.Item(4, 3).Value = '15460000'
$ws.Cells.Item(4, 4).Value = '17530612'
$ws.Cells.Item(4, 5).Value = '1546-07-02'
$ws.Cells.Item(4, 6).Value = '1753-06-12'
$ws.Cells.Item(4, 7).Value = '1546'
$ws.Cells.Item(4, 8).Value = '1753'
$ws.Cells.Item(5, 1).Value = 'Paris'
$ws.Cells.Item(5, 2).Value = 44
$ws.Cells.Item(5, 3).Value = '15340815'
$ws.Cells.Item(5, 4).Value = '17590310'
$ws.Cells.Item(5, 5).Value = '1534-08-15'
$ws.Cells.Item(5, 6).Value = '1759-03-10'
$ws.Cells.Item(5, 7).Value = '1534'
$ws.Cells.Item(5, 8).Value = '1759'
$ws.Cells.Item(6, 1).Value = 'Roma'
$ws.Cells.Item(6, 2).Value = 41
$ws.Cells.Item(6, 3).Value = '15400927'
$ws.Cells.Item(6, 4).Value = '17560709'
$ws.Cells.Item(6, 5).Value = '1540-09-27'
$ws.Cells.Item(6, 6).Value = '1756-07-09'
$ws.Cells.Item(6, 7).Value = '1540'
$ws.Cells.Item(6, 8).Value = '1756'
$ws.Cells.Item(7, 1).Value = 'Goa'
$ws.Cells.Item(7, 2).Value = 33
$ws.Cells.Item(7, 3).Value = '15480000'
$ws.Cells.Item(7, 4).Value = '17360000'
$ws.Cells.Item(7, 5).Value = '1548-07-02'
$ws.Cells.Item(7, 6).Value = '1736-07-02'
$ws.Cells.Item(7, 7).Value = '1548'
$ws.Cells.Item(7, 8).Value = '1736'
$ws.Cells.Item(8, 1).Value = 'Évora'
$ws.Cells.Item(8, 2).Value = 26
$ws.Cells.Item(8, 3).Value = '15660714'
$ws.Cells.Item(8, 4).Value = '17460405'
$ws.Cells.Item(8, 5).Value = '1566-07-14'
$ws.Cells.Item(8, 6).Value = '1746-04-05'
$ws.Cells.Item(8, 7).Value = '1566'
$ws.Cells.Item(8, 8).Value = '1746'
$ws.Cells.Item(9, 1).Value = 'Nancy'
$ws.Cells.Item(9, 2).Value = 15
$ws.Cells.Item(9, 3).Value = '16270929'
$ws.Cells.Item(9, 4).Value = '17510827'
$ws.Cells.Item(9, 5).Value = '1627-09-29'
$ws.Cells.Item(9, 6).Value = '1751-08-27'
$ws.Cells.Item(9, 7).Value = '1627'
$ws.Cells.Item(9, 8).Value = '1751'
$ws.Cells.Item(10, 1).Value = 'Macau'
$ws.Cells.Item(10, 2).Value = 14
$ws.Cells.Item(10, 3).Value = '16280000'
$ws.Cells.Item(10, 4).Value = '17490201'
$ws.Cells.Item(10, 5).Value = '1628-07-02'
$ws.Cells.Item(10, 6).Value = '1749-02-01'
$ws.Cells.Item(10, 7).Value = '1628'
$ws.Cells.Item(10, 8).Value = '1749'
$ws.Cells.Item(11, 1).Value = 'Avignon'
$ws.Cells.Item(11, 2).Value = 13
$ws.Cells.Item(11, 3).Value = '16150926'
$ws.Cells.Item(11, 4).Value = '17370927'
$ws.Cells.Item(11, 5).Value = '1615-09-26'
$ws.Cells.Item(11, 6).Value = '1737-09-27'
$ws.Cells.Item(11, 7).Value = '1615'
$ws.Cells.Item(11, 8).Value = '1737'
$ws.Cells.Item(12, 1).Value = 'Landsberg'
$ws.Cells.Item(12, 2).Value = 13
$ws.Cells.Item(12, 3).Value = '16230729'
$ws.Cells.Item(12, 4).Value = '17571009'
$ws.Cells.Item(12, 5).Value = '1623-07-29'
$ws.Cells.Item(12, 6).Value = '1757-10-09'
$ws.Cells.Item(12, 7).Value = '1623'
$ws.Cells.Item(12, 8).Value = '1757'
$ws.Cells.Item(13, 1).Value = 'Viena'
$ws.Cells.Item(13, 2).Value = 13
$ws.Cells.Item(13, 3).Value = '16271016'
$ws.Cells.Item(13, 4).Value = '17651018'
$ws.Cells.Item(13, 5).Value = '1627-10-16'
$ws.Cells.Item(13, 6).Value = '1765-10-18'
$ws.Cells.Item(13, 7).Value = '1627'
$ws.Cells.Item(13, 8).Value = '1765'
$ws.Cells.Item(14, 1).Value = 'Pequim'
$ws.Cells.Item(14, 2).Value = 11
$ws.Cells.Item(14, 3).Value = '17310308'
$ws.Cells.Item(14, 4).Value = '17730214'
$ws.Cells.Item(14, 5).Value = '1731-03-08'
$ws.Cells.Item(14, 6).Value = '1773-02-14'
$ws.Cells.Item(14, 7).Value = '1731'
$ws.Cells.Item(14, 8).Value = '1773'
$ws.Cells.Item(15, 1).Value = 'Mechelen'
$ws.Cells.Item(15, 2).Value = 10
$ws.Cells.Item(15, 3).Value = '16341015'
$ws.Cells.Item(15, 4).Value = '16720424'
$ws.Cells.Item(15, 5).Value = '1634-10-15'
$ws.Cells.Item(15, 6).Value = '1672-04-24'
$ws.Cells.Item(15, 7).Value = '1634'
$ws.Cells.Item(15, 8).Value = '1672'
$ws.Cells.Item(16, 1).Value = 'Nápoles'
$ws.Cells.Item(16, 2).Value = 10
$ws.Cells.Item(16, 3).Value = '15560120'
$ws.Cells.Item(16, 4).Value = '17000311'
$ws.Cells.Item(16, 5).Value = '1556-01-20'
$ws.Cells.Item(16, 6).Value = '1700-03-11'
$ws.Cells.Item(16, 7).Value = '1556'
$ws.Cells.Item(16, 8).Value = '1700'
$ws.Cells.Item(17, 1).Value = 'Génova'
$ws.Cells.Item(17, 2).Value = 9
$ws.Cells.Item(17, 3).Value = '16390915'
$ws.Cells.Item(17, 4).Value = '17680930'
$ws.Cells.Item(17, 5).Value = '1639-09-15'
$ws.Cells.Item(17, 6).Value = '1768-09-30'
$ws.Cells.Item(17, 7).Value = '1639'
$ws.Cells.Item(17, 8).Value = '1768'
$ws.Cells.Item(18, 1).Value = 'Tournai'
$ws.Cells.Item(18, 2).Value = 9
$ws.Cells.Item(18, 3).Value = '15960705'
$ws.Cells.Item(18, 4).Value = '16960929'
$ws.Cells.Item(18, 5).Value = '1596-07-05'
$ws.Cells.Item(18, 6).Value = '1696-09-29'
$ws.Cells.Item(18, 7).Value = '1596'
$ws.Cells.Item(18, 8).Value = '1696'
$ws.Cells.Item(19, 1).Value = 'Bordeaux'
$ws.Cells.Item(19, 2).Value = 8
$ws.Cells.Item(19, 3).Value = '16680920'
$ws.Cells.Item(19, 4).Value = '17431107'
$ws.Cells.Item(19, 5).Value = '1668-09-20'
$ws.Cells.Item(19, 6).Value = '1743-11-07'
$ws.Cells.Item(19, 7).Value = '1668'
$ws.Cells.Item(19, 8).Value = '1743'
$ws.Cells.Item(20, 1).Value = 'Alcalá'
$ws.Cells.Item(20, 2).Value = 6
$ws.Cells.Item(20, 3).Value = '15550925'
$ws.Cells.Item(20, 4).Value = '15710624'
$ws.Cells.Item(20, 5).Value = '1555-09-25'
$ws.Cells.Item(20, 6).Value = '1571-06-24'
$ws.Cells.Item(20, 7).Value = '1555'
$ws.Cells.Item(20, 8).Value = '1571'
$ws.Cells.Item(21, 1).Value = 'Brno'
$ws.Cells.Item(21, 2).Value = 6
$ws.Cells.Item(21, 3).Value = '16060000'
$ws.Cells.Item(21, 4).Value = '17261009'
$ws.Cells.Item(21, 5).Value = '1606-07-02'
$ws.Cells.Item(21, 6).Value = '1726-10-09'
$ws.Cells.Item(21, 7).Value = '1606'
$ws.Cells.Item(21, 8).Value = '1726'
$ws.Cells.Item(22, 1).Value = 'Chieri'
$ws.Cells.Item(22, 2).Value = 5
$ws.Cells.Item(22, 3).Value = '16280213'
$ws.Cells.Item(22, 4).Value = '16880120'
$ws.Cells.Item(22, 5).Value = '1628-02-13'
$ws.Cells.Item(22, 6).Value = '1688-01-20'
$ws.Cells.Item(22, 7).Value = '1628'
$ws.Cells.Item(22, 8).Value = '1688'
$ws.Cells.Item(23, 1).Value = 'Japão'
$ws.Cells.Item(23, 2).Value = 5
$ws.Cells.Item(23, 3).Value = '15560000'
$ws.Cells.Item(23, 4).Value = '16260600'
$ws.Cells.Item(23, 5).Value = '1556-07-02'
$ws.Cells.Item(23, 6).Value = '1626-06-15'
$ws.Cells.Item(23, 7).Value = '1556'
$ws.Cells.Item(23, 8).Value = '1626'
$ws.Cells.Item(24, 1).Value = 'Cracóvia'
$ws.Cells.Item(24, 2).Value = 4
$ws.Cells.Item(24, 3).Value = '16310816'
$ws.Cells.Item(24, 4).Value = '16760909'
$ws.Cells.Item(24, 5).Value = '1631-08-16'
$ws.Cells.Item(24, 6).Value = '1676-09-09'
$ws.Cells.Item(24, 7).Value = '1631'
$ws.Cells.Item(24, 8).Value = '1676'
$ws.Cells.Item(25, 1).Value = 'Lyon'
$ws.Cells.Item(25, 2).Value = 4
$ws.Cells.Item(25, 3).Value = '16740927'
$ws.Cells.Item(25, 4).Value = '16930923'
$ws.Cells.Item(25, 5).Value = '1674-09-27'
$ws.Cells.Item(25, 6).Value = '1693-09-23'
$ws.Cells.Item(25, 7).Value = '1674'
$ws.Cells.Item(25, 8).Value = '1693'
$ws.Cells.Item(26, 1).Value = 'Palermo'
$ws.Cells.Item(26, 2).Value = 4
$ws.Cells.Item(26, 3).Value = '16181103'
$ws.Cells.Item(26, 4).Value = '16561022'
$ws.Cells.Item(26, 5).Value = '1618-11-03'
$ws.Cells.Item(26, 6).Value = '1656-10-22'
$ws.Cells.Item(26, 7).Value = '1618'
$ws.Cells.Item(26, 8).Value = '1656'
$ws.Cells.Item(27, 1).Value = 'Portugal'
$ws.Cells.Item(27, 2).Value = 4
$ws.Cells.Item(27, 3).Value = '15480613'
$ws.Cells.Item(27, 4).Value = '16780000'
$ws.Cells.Item(27, 5).Value = '1548-06-13'
$ws.Cells.Item(27, 6).Value = '1678-07-02'
$ws.Cells.Item(27, 7).Value = '1548'
$ws.Cells.Item(27, 8).Value = '1678'
$ws.Cells.Item(28, 1).Value = 'Shiuchow'
$ws.Cells.Item(28, 2).Value = 4
$ws.Cells.Item(28, 3).Value = '15891100'
$ws.Cells.Item(28, 4).Value = '16050000'
$ws.Cells.Item(28, 5).Value = '1589-11-15'
$ws.Cells.Item(28, 6).Value = '1605-07-02'
$ws.Cells.Item(28, 7).Value = '1589'
$ws.Cells.Item(28, 8).Value = '1605'
$ws.Cells.Item(29, 1).Value = 'Bolonha'
$ws.Cells.Item(29, 2).Value = 3
$ws.Cells.Item(29, 3).Value = '16680000'
$ws.Cells.Item(29, 4).Value = '17250728'
$ws.Cells.Item(29, 5).Value = '1668-07-02'
$ws.Cells.Item(29, 6).Value = '1725-07-28'
$ws.Cells.Item(29, 7).Value = '1668'
$ws.Cells.Item(29, 8).Value = '1725'
$ws.Cells.Item(30, 1).Value = 'Hangchow'
$ws.Cells.Item(30, 2).Value = 3
$ws.Cells.Item(30, 3).Value = '16270000'
$ws.Cells.Item(30, 4).Value = '16790621'
$ws.Cells.Item(30, 5).Value = '1627-07-02'
$ws.Cells.Item(30, 6).Value = '1679-06-21'
$ws.Cells.Item(30, 7).Value = '1627'
$ws.Cells.Item(30, 8).Value = '1679'
$ws.Cells.Item(31, 1).Value = 'Milão'
$ws.Cells.Item(31, 2).Value = 3
$ws.Cells.Item(31, 3).Value = '16731021'
$ws.Cells.Item(31, 4).Value = '16901101'
$ws.Cells.Item(31, 5).Value = '1673-10-21'
$ws.Cells.Item(31, 6).Value = '1690-11-01'
$ws.Cells.Item(31, 7).Value = '1673'
$ws.Cells.Item(31, 8).Value = '1690'
$ws.Cells.Item(32, 1).Value = 'Toulouse'
$ws.Cells.Item(32, 2).Value = 3
$ws.Cells.Item(32, 3).Value = '16800921'
$ws.Cells.Item(32, 4).Value = '17280930'
$ws.Cells.Item(32, 5).Value = '1680-09-21'
$ws.Cells.Item(32, 6).Value = '1728-09-30'
$ws.Cells.Item(32, 7).Value = '1680'
$ws.Cells.Item(32, 8).Value = '1728'
$ws.Cells.Item(33, 1).Value = 'Trier'
$ws.Cells.Item(33, 2).Value = 3
$ws.Cells.Item(33, 3).Value = '16770821'
$ws.Cells.Item(33, 4).Value = '17271019'
$ws.Cells.Item(33, 5).Value = '1677-08-21'
$ws.Cells.Item(33, 6).Value = '1727-10-19'
$ws.Cells.Item(33, 7).Value = '1677'
$ws.Cells.Item(33, 8).Value = '1727'
$ws.Cells.Item(34, 1).Value = 'Arona'
$ws.Cells.Item(34, 2).Value = 2
$ws.Cells.Item(34, 3).Value = '15960921'
$ws.Cells.Item(34, 4).Value = '16140824'
$ws.Cells.Item(34, 5).Value = '1596-09-21'
$ws.Cells.Item(34, 6).Value = '1614-08-24'
$ws.Cells.Item(34, 7).Value = '1596'
$ws.Cells.Item(34, 8).Value = '1614'
$ws.Cells.Item(35, 1).Value = 'Boémia'
$ws.Cells.Item(35, 2).Value = 2
$ws.Cells.Item(35, 3).Value = '17290927'
$ws.Cells.Item(35, 4).Value = '17291010'
$ws.Cells.Item(35, 5).Value = '1729-09-27'
$ws.Cells.Item(35, 6).Value = '1729-10-10'
$ws.Cells.Item(35, 7).Value = '1729'
$ws.Cells.Item(35, 8).Value = '1729'
$ws.Cells.Item(36, 1).Value = 'Douai'
$ws.Cells.Item(36, 2).Value = 2
$ws.Cells.Item(36, 3).Value = '15941109'
$ws.Cells.Item(36, 4).Value = '15990512'
$ws.Cells.Item(36, 5).Value = '1594-11-09'
$ws.Cells.Item(36, 6).Value = '1599-05-12'
$ws.Cells.Item(36, 7).Value = '1594'
$ws.Cells.Item(36, 8).Value = '1599'
$ws.Cells.Item(37, 1).Value = 'Ferrara'
$ws.Cells.Item(37, 2).Value = 2
$ws.Cells.Item(37, 3).Value = '15550000'
$ws.Cells.Item(37, 4).Value = '15561200'
$ws.Cells.Item(37, 5).Value = '1555-07-02'
$ws.Cells.Item(37, 6).Value = '1556-12-15'
$ws.Cells.Item(37, 7).Value = '1555'
$ws.Cells.Item(37, 8).Value = '1556'
$ws.Cells.Item(38, 1).Value = 'Japão (província)'
$ws.Cells.Item(38, 2).Value = 2
$ws.Cells.Item(38, 3).Value = '16910113'
$ws.Cells.Item(38, 4).Value = '17280523'
$ws.Cells.Item(38, 5).Value = '1691-01-13'
$ws.Cells.Item(38, 6).Value = '1728-05-23'
$ws.Cells.Item(38, 7).Value = '1691'
$ws.Cells.Item(38, 8).Value = '1728'
$ws.Cells.Item(39, 1).Value = 'Leoben'
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = '16871224'
$ws.Cells.Item(39, 4).Value = '16881012'
$ws.Cells.Item(39, 5).Value = '1687-12-24'
$ws.Cells.Item(39, 6).Value = '1688-10-12'
$ws.Cells.Item(39, 7).Value = '1687'
$ws.Cells.Item(39, 8).Value = '1688'
$ws.Cells.Item(40, 1).Value = 'Messina'
$ws.Cells.Item(40, 2).Value = 2
$ws.Cells.Item(40, 3).Value = '15820000'
$ws.Cells.Item(40, 4).Value = '17061207'
$ws.Cells.Item(40, 5).Value = '1582-07-02'
$ws.Cells.Item(40, 6).Value = '1706-12-07'
$ws.Cells.Item(40, 7).Value = '1582'
$ws.Cells.Item(40, 8).Value = '1706'
$ws.Cells.Item(41, 1).Value = 'Novellara'
$ws.Cells.Item(41, 2).Value = 2
$ws.Cells.Item(41, 3).Value = '16001101'
$ws.Cells.Item(41, 4).Value = '16170121'
$ws.Cells.Item(41, 5).Value = '1600-11-01'
$ws.Cells.Item(41, 6).Value = '1617-01-21'
$ws.Cells.Item(41, 7).Value = '1600'
$ws.Cells.Item(41, 8).Value = '1617'
$ws.Cells.Item(42, 1).Value = 'Saragoça'
$ws.Cells.Item(42, 2).Value = 2
$ws.Cells.Item(42, 3).Value = '16740613'
$ws.Cells.Item(42, 4).Value = '16771115'
$ws.Cells.Item(42, 5).Value = '1674-06-13'
$ws.Cells.Item(42, 6).Value = '1677-11-15'
$ws.Cells.Item(42, 7).Value = '1674'
$ws.Cells.Item(42, 8).Value = '1677'
$ws.Cells.Item(43, 1).Value = 'Vilnius'
$ws.Cells.Item(43, 2).Value = 2
$ws.Cells.Item(43, 3).Value = '16180531'
$ws.Cells.Item(43, 4).Value = '16720811'
$ws.Cells.Item(43, 5).Value = '1618-05-31'
$ws.Cells.Item(43, 6).Value = '1672-08-11'
$ws.Cells.Item(43, 7).Value = '1618'
$ws.Cells.Item(43, 8).Value = '1672'
$ws.Cells.Item(44, 1).Value = 'Bahia'
$ws.Cells.Item(44, 2).Value = 1
$ws.Cells.Item(44, 3).Value = '16410000'
$ws.Cells.Item(44, 4).Value = '16410000'
$ws.Cells.Item(44, 5).Value = '1641-07-02'
$ws.Cells.Item(44, 6).Value = '1641-07-02'
$ws.Cells.Item(44, 7).Value = '1641'
$ws.Cells.Item(44, 8).Value = '1641'
$ws.Cells.Item(45, 1).Value = 'China'
$ws.Cells.Item(45, 2).Value = 1
$ws.Cells.Item(45, 3).Value = '17280105'
$ws.Cells.Item(45, 4).Value = '17280105'
$ws.Cells.Item(45, 5).Value = '1728-01-05'
$ws.Cells.Item(45, 6).Value = '1728-01-05'
$ws.Cells.Item(45, 7).Value = '1728'
$ws.Cells.Item(45, 8).Value = '1728'
$ws.Cells.Item(46, 1).Value = 'Colorno'
$ws.Cells.Item(46, 2).Value = 1
$ws.Cells.Item(46, 3).Value = '17991116'
$ws.Cells.Item(46, 4).Value = '17991116'
$ws.Cells.Item(46, 5).Value = '1799-11-16'
$ws.Cells.Item(46, 6).Value = '1799-11-16'
$ws.Cells.Item(46, 7).Value = '1799'
$ws.Cells.Item(46, 8).Value = '1799'
$ws.Cells.Item(47, 1).Value = 'Courtrai'
$ws.Cells.Item(47, 2).Value = 1
$ws.Cells.Item(47, 3).Value = '16440926'
$ws.Cells.Item(47, 4).Value = '16440926'
$ws.Cells.Item(47, 5).Value = '1644-09-26'
$ws.Cells.Item(47, 6).Value = '1644-09-26'
$ws.Cells.Item(47, 7).Value = '1644'
$ws.Cells.Item(47, 8).Value = '1644'
$ws.Cells.Item(48, 1).Value = 'Krems'
$ws.Cells.Item(48, 2).Value = 1
$ws.Cells.Item(48, 3).Value = '16641031'
$ws.Cells.Item(48, 4).Value = '16641031'
$ws.Cells.Item(48, 5).Value = '1664-10-31'
$ws.Cells.Item(48, 6).Value = '1664-10-31'
$ws.Cells.Item(48, 7).Value = '1664'
$ws.Cells.Item(48, 8).Value = '1664'
$ws.Cells.Item(49, 1).Value = 'Lima, Peru'
$ws.Cells.Item(49, 2).Value = 1
$ws.Cells.Item(49, 3).Value = '15680711'
$ws.Cells.Item(49, 4).Value = '15680711'
$ws.Cells.Item(49, 5).Value = '1568-07-11'
$ws.Cells.Item(49, 6).Value = '1568-07-11'
$ws.Cells.Item(49, 7).Value = '1568'
$ws.Cells.Item(49, 8).Value = '1568'
$ws.Cells.Item(50, 1).Value = 'Lorette'
$ws.Cells.Item(50, 2).Value = 1
$ws.Cells.Item(50, 3).Value = '15590425'
$ws.Cells.Item(50, 4).Value = '15590425'
$ws.Cells.Item(50, 5).Value = '1559-04-25'
$ws.Cells.Item(50, 6).Value = '1559-04-25'
$ws.Cells.Item(50, 7).Value = '1559'
$ws.Cells.Item(50, 8).Value = '1559'
$ws.Cells.Item(51, 1).Value = 'Mainz'
$ws.Cells.Item(51, 2).Value = 1
$ws.Cells.Item(51, 3).Value = '16730717'
$ws.Cells.Item(51, 4).Value = '16730717'
$ws.Cells.Item(51, 5).Value = '1673-07-17'
$ws.Cells.Item(51, 6).Value = '1673-07-17'
$ws.Cells.Item(51, 7).Value = '1673'
$ws.Cells.Item(51, 8).Value = '1673'
$ws.Cells.Item(52, 1).Value = 'Manila'
$ws.Cells.Item(52, 2).Value = 1
$ws.Cells.Item(52, 3).Value = '16721011'
$ws.Cells.Item(52, 4).Value = '16721011'
$ws.Cells.Item(52, 5).Value = '1672-10-11'
$ws.Cells.Item(52, 6).Value = '1672-10-11'
$ws.Cells.Item(52, 7).Value = '1672'
$ws.Cells.Item(52, 8).Value = '1672'
$ws.Cells.Item(53, 1).Value = 'Mazowsze (província)'
$ws.Cells.Item(53, 2).Value = 1
$ws.Cells.Item(53, 3).Value = '17700813'
$ws.Cells.Item(53, 4).Value = '17700813'
$ws.Cells.Item(53, 5).Value = '1770-08-13'
$ws.Cells.Item(53, 6).Value = '1770-08-13'
$ws.Cells.Item(53, 7).Value = '1770'
$ws.Cells.Item(53, 8).Value = '1770'
$ws.Cells.Item(54, 1).Value = 'Milão (província)'
$ws.Cells.Item(54, 2).Value = 1
$ws.Cells.Item(54, 3).Value = '16581105'
$ws.Cells.Item(54, 4).Value = '16581105'
$ws.Cells.Item(54, 5).Value = '1658-11-05'
$ws.Cells.Item(54, 6).Value = '1658-11-05'
$ws.Cells.Item(54, 7).Value = '1658'
$ws.Cells.Item(54, 8).Value = '1658'
$ws.Cells.Item(55, 1).Value = 'Nan-tch''ang'
$ws.Cells.Item(55, 2).Value = 1
$ws.Cells.Item(55, 3).Value = '16080000'
$ws.Cells.Item(55, 4).Value = '16080000'
$ws.Cells.Item(55, 5).Value = '1608-07-02'
$ws.Cells.Item(55, 6).Value = '1608-07-02'
$ws.Cells.Item(55, 7).Value = '1608'
$ws.Cells.Item(55, 8).Value = '1608'
$ws.Cells.Item(56, 1).Value = 'Nanquim'
$ws.Cells.Item(56, 2).Value = 1
$ws.Cells.Item(56, 3).Value = '16080300'
$ws.Cells.Item(56, 4).Value = '16080300'
$ws.Cells.Item(56, 5).Value = '1608-03-15'
$ws.Cells.Item(56, 6).Value = '1608-03-15'
$ws.Cells.Item(56, 7).Value = '1608'
$ws.Cells.Item(56, 8).Value = '1608'
$ws.Cells.Item(57, 1).Value = 'Ormuz'
$ws.Cells.Item(57, 2).Value = 1
$ws.Cells.Item(57, 3).Value = '000000'
$ws.Cells.Item(57, 4).Value = '000000'
$ws.Cells.Item(57, 5).Value = $null
$ws.Cells.Item(57, 6).Value = $null
$ws.Cells.Item(57, 7).Value = $null
$ws.Cells.Item(57, 8).Value = $null
$ws.Cells.Item(58, 1).Value = 'Ozukio (noviciado)'
$ws.Cells.Item(58, 2).Value = 1
$ws.Cells.Item(58, 3).Value = '15811100'
$ws.Cells.Item(58, 4).Value = '15811100'
$ws.Cells.Item(58, 5).Value = '1581-11-15'
$ws.Cells.Item(58, 6).Value = '1581-11-15'
$ws.Cells.Item(58, 7).Value = '1581'
$ws.Cells.Item(58, 8).Value = '1581'
$ws.Cells.Item(59, 1).Value = 'Polotsk'
$ws.Cells.Item(59, 2).Value = 1
$ws.Cells.Item(59, 3).Value = '17860903'
$ws.Cells.Item(59, 4).Value = '17860903'
$ws.Cells.Item(59, 5).Value = '1786-09-03'
$ws.Cells.Item(59, 6).Value = '1786-09-03'
$ws.Cells.Item(59, 7).Value = '1786'
$ws.Cells.Item(59, 8).Value = '1786'
$ws.Cells.Item(60, 1).Value = 'Salamanca'
$ws.Cells.Item(60, 2).Value = 1
$ws.Cells.Item(60, 3).Value = '16830419'
$ws.Cells.Item(60, 4).Value = '16830419'
$ws.Cells.Item(60, 5).Value = '1683-04-19'
$ws.Cells.Item(60, 6).Value = '1683-04-19'
$ws.Cells.Item(60, 7).Value = '1683'
$ws.Cells.Item(60, 8).Value = '1683'
$ws.Cells.Item(61, 1).Value = 'Shanghai'
$ws.Cells.Item(61, 2).Value = 1
$ws.Cells.Item(61, 3).Value = '16100000'
$ws.Cells.Item(61, 4).Value = '16100000'
$ws.Cells.Item(61, 5).Value = '1610-07-02'
$ws.Cells.Item(61, 6).Value = '1610-07-02'
$ws.Cells.Item(61, 7).Value = '1610'
$ws.Cells.Item(61, 8).Value = '1610'
$ws.Cells.Item(62, 1).Value = 'Todos-os-Santos, Nagasaki'
$ws.Cells.Item(62, 2).Value = 1
$ws.Cells.Item(62, 3).Value = '16070202'
$ws.Cells.Item(62, 4).Value = '16070202'
$ws.Cells.Item(62, 5).Value = '1607-02-02'
$ws.Cells.Item(62, 6).Value = '1607-02-02'
$ws.Cells.Item(62, 7).Value = '1607'
$ws.Cells.Item(62, 8).Value = '1607'
$ws.Cells.Item(63, 1).Value = 'Toulouse (província)'
$ws.Cells.Item(63, 2).Value = 1
$ws.Cells.Item(63, 3).Value = '17500319'
$ws.Cells.Item(63, 4).Value = '17500319'
$ws.Cells.Item(63, 5).Value = '1750-03-19'
$ws.Cells.Item(63, 6).Value = '1750-03-19'
$ws.Cells.Item(63, 7).Value = '1750'
$ws.Cells.Item(63, 8).Value = '1750'
$ws.Cells.Item(64, 1).Value = 'Trenčín'
$ws.Cells.Item(64, 2).Value = 1
$ws.Cells.Item(64, 3).Value = '17291027'
$ws.Cells.Item(64, 4).Value = '17291027'
$ws.Cells.Item(64, 5).Value = '1729-10-27'
$ws.Cells.Item(64, 6).Value = '1729-10-27'
$ws.Cells.Item(64, 7).Value = '1729'
$ws.Cells.Item(64, 8).Value = '1729'
$ws.Cells.Item(65, 1).Value = 'Valença'
$ws.Cells.Item(65, 2).Value = 1
$ws.Cells.Item(65, 3).Value = '15610927'
$ws.Cells.Item(65, 4).Value = '15610927'
$ws.Cells.Item(65, 5).Value = '1561-09-27'
$ws.Cells.Item(65, 6).Value = '1561-09-27'
$ws.Cells.Item(65, 7).Value = '1561'
$ws.Cells.Item(65, 8).Value = '1561'
$ws.Cells.Item(66, 1).Value = 'Veneza'
$ws.Cells.Item(66, 2).Value = 1
$ws.Cells.Item(66, 3).Value = '17180424'
$ws.Cells.Item(66, 4).Value = '17180424'
$ws.Cells.Item(66, 5).Value = '1718-04-24'
$ws.Cells.Item(66, 6).Value = '1718-04-24'
$ws.Cells.Item(66, 7).Value = '1718'
$ws.Cells.Item(66, 8).Value = '1718'
$ws.Cells.Item(67, 1).Value = 'Villaregio'
$ws.Cells.Item(67, 2).Value = 1
$ws.Cells.Item(67, 3).Value = '15890406'
$ws.Cells.Item(67, 4).Value = '15890406'
$ws.Cells.Item(67, 5).Value = '1589-04-06'
$ws.Cells.Item(67, 6).Value = '1589-04-06'
$ws.Cells.Item(67, 7).Value = '1589'
$ws.Cells.Item(67, 8).Value = '1589'
$ws.Cells.Item(68, 1).Value = 'Índia'
$ws.Cells.Item(68, 2).Value = 1
$ws.Cells.Item(68, 3).Value = '15530000'
$ws.Cells.Item(68, 4).Value = '15530000'
$ws.Cells.Item(68, 5).Value = '1553-07-02'
$ws.Cells.Item(68, 6).Value = '1553-07-02'
$ws.Cells.Item(68, 7).Value = '1553'
$ws.Cells.Item(68, 8).Value = '1553'
